$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column to the left (M) so the newly
# inserted column can inherit the same custom width, just like Excel
# does when a column is inserted.
$leftWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new (blank) column before column N; this pushes the old
# N, O, P columns one to the right (-> O, P, Q) together with their
# values/styles.
$ws.Columns("N:N").Insert()

# The newly inserted column keeps the neighbour's width.
$ws.Columns("N:N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab and put the
# selection on J20.
$ws.Activate() | Out-Null
$ws.Range("J20").Select() | Out-Null
